$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.020.89"
$ws.Range("E2").Value = "  -2.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.068.29"
$ws.Range("E3").Value = "  -1.69%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "520.88"
$ws.Range("E5").Value = "  -1.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.36"
$ws.Range("E6").Value = "  -4.93%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.067.84"
$ws.Range("E8").Value = "  -1.60%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.458"
$ws.Range("E9").Value = "  +3.02%  "
$ws.Range("E10").Value = "  +2.02%  "
$ws.Range("E11").Value = "  -2.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.398"
$ws.Range("E12").Value = "  +1.22%  "
$ws.Range("E13").Value = "  +0.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.597.17"
$ws.Range("E14").Value = "  -1.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.19"
$ws.Range("E15").Value = "  -2.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "57.039.43"
$ws.Range("E17").Value = "  -1.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.065.55"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.86"
$ws.Range("E19").Value = "  -4.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.42"
$ws.Range("E20").Value = "  -2.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.81"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "349.06"
$ws.Range("E22").Value = "  +1.54%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.71"
$ws.Range("E24").Value = "  +1.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.497"
$ws.Range("E25").Value = "  -3.66%  "
$ws.Range("E26").Value = "  -2.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("E28").Value = "  -7.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.18"
$ws.Range("E30").Value = "  -1.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.86"
$ws.Range("E31").Value = "  -1.47%  "
$ws.Range("E32").Value = "  -8.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.86"
$ws.Range("E33").Value = "  -1.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "159.20"
$ws.Range("E34").Value = "  +0.48%  "
$ws.Range("E35").Value = "  +3.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.12"
$ws.Range("E36").Value = "  -5.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.98"
$ws.Range("E37").Value = "  -3.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.34"
$ws.Range("E38").Value = "  -4.35%  "
$ws.Range("E39").Value = "  -2.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0656"
$ws.Range("E40").Value = "  -2.37%  "
$ws.Range("E41").Value = "  -4.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.04"
$ws.Range("E42").Value = "  +0.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.690"
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.389.72"
$ws.Range("E44").Value = "  +5.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "36.62"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.107.00"
$ws.Range("E47").Value = "  -1.64%  "
$ws.Range("E48").Value = "  -1.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.949"
$ws.Range("E49").Value = "  -5.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.95"
$ws.Range("E50").Value = "  -2.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.60"
$ws.Range("E51").Value = "  -5.18%  "
